$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $text) {
    $rng = $ws.Range($cellAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "26.976.39"
Set-TextValue "E2" "  -0.33%  "
Set-TextValue "D3" "1.823.42"
Set-TextValue "E3" "  +0.03%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  -0.84%  "
Set-TextValue "D5" "311.32"
Set-TextValue "E5" "  -0.06%  "
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  -0.69%  "
Set-TextValue "E7" "  -0.72%  "
Set-TextValue "D8" "0.3715"
Set-TextValue "E8" "  +2.10%  "
Set-TextValue "D9" "0.07336"
Set-TextValue "E9" "  +0.53%  "
Set-TextValue "D10" "0.8756"
Set-TextValue "E10" "  +0.94%  "
Set-TextValue "D11" "0.07867"
Set-TextValue "E11" "  +3.27%  "
Set-TextValue "D12" "19.72"
Set-TextValue "E12" "  -2.21%  "
Set-TextValue "D13" "1.835.54"
Set-TextValue "E13" "  -0.18%  "
Set-TextValue "D14" "5.336"
Set-TextValue "E14" "  -0.11%  "
Set-TextValue "D15" "6.554"
Set-TextValue "E15" "  +1.25%  "
Set-TextValue "D16" "91.20"
Set-TextValue "E16" "  -1.74%  "
Set-TextValue "E17" "  -0.47%  "
Set-TextValue "D18" "0.000008829"
Set-TextValue "E18" "  +2.20%  "
Set-TextValue "D20" "14.81"
Set-TextValue "E20" "  +2.14%  "
Set-TextValue "D21" "26.988.85"
Set-TextValue "E21" "  -1.77%  "
Set-TextValue "D22" "5.091"
Set-TextValue "E22" "  -2.01%  "
Set-TextValue "E23" "  -0.44%  "
Set-TextValue "D24" "2.098.06"
Set-TextValue "E24" "  +0.39%  "
Set-TextValue "D25" "153.11"
Set-TextValue "E25" "  +0.87%  "
Set-TextValue "D26" "1.848"
Set-TextValue "E26" "  -0.90%  "
Set-TextValue "D27" "18.42"
Set-TextValue "E27" "  +0.75%  "
Set-TextValue "D28" "2.034"
Set-TextValue "E28" "  -3.27%  "
Set-TextValue "D29" "5.124"
Set-TextValue "D30" "115.62"
Set-TextValue "E30" "  -0.48%  "
Set-TextValue "D31" "0.08869"
Set-TextValue "E31" "  -0.67%  "
Set-TextValue "D32" "2.952"
Set-TextValue "E32" "  -0.19%  "
Set-TextValue "D33" "0.7287"
Set-TextValue "E33" "  -0.43%  "
Set-TextValue "D34" "4.430"
Set-TextValue "E34" "  -0.53%  "
Set-TextValue "D35" "1.130"
Set-TextValue "E35" "  -0.89%  "
Set-TextValue "D36" "2.474"
Set-TextValue "E36" "  -2.05%  "
Set-TextValue "D37" "0.01946"
Set-TextValue "E37" "  +1.42%  "
Set-TextValue "D38" "1.066"
Set-TextValue "E38" "  -0.69%  "
Set-TextValue "D39" "0.05213"
Set-TextValue "E39" "  -1.06%  "
Set-TextValue "D40" "2.946"
Set-TextValue "E40" "  +0.15%  "
Set-TextValue "D41" "7.089"
Set-TextValue "E41" "  -0.48%  "
Set-TextValue "D42" "0.5148"
Set-TextValue "E42" "  -1.30%  "
Set-TextValue "D43" "0.1620"
Set-TextValue "E43" "  -0.82%  "
Set-TextValue "D44" "8.155"
Set-TextValue "E44" "  -1.24%  "
Set-TextValue "D45" "0.4828"
Set-TextValue "E45" "  -0.63%  "
Set-TextValue "E46" "  +0.83%  "
Set-TextValue "D47" "1.002"
Set-TextValue "E47" "  -0.75%  "
Set-TextValue "D48" "102.43"
Set-TextValue "E48" "  -1.28%  "
Set-TextValue "D49" "1.626"
Set-TextValue "E49" "  -1.02%  "
Set-TextValue "D50" "0.06193"
Set-TextValue "D51" "64.74"
Set-TextValue "E51" "  +0.21%  "
